$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET ---
$taskSheet = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Name of the person the sheet belongs to
$taskSheet.Range("C1").Value = "Jesse Hare"

# Task rows: Stage / Task / Estimated Work Remaining at Start of Week / Hours Spent this Week / New Estimate
$taskSheet.Range("A3").Value = "Project Build"
$taskSheet.Range("B3").Value = "Finalise Search algorithm"
$taskSheet.Range("C3").Value = 15
$taskSheet.Range("D3").Value = 9
$taskSheet.Range("E3").Value = 0

$taskSheet.Range("A4").Value = "Project Build"
$taskSheet.Range("B4").Value = "Work on sorting algorithm"
$taskSheet.Range("C4").Value = 15
$taskSheet.Range("D4").Value = 6
$taskSheet.Range("E4").Value = 0

$taskSheet.Range("A5").Value = "Project Build"
$taskSheet.Range("B5").Value = "Re-design GUI elements"
$taskSheet.Range("C5").Value = 3
$taskSheet.Range("D5").Value = 4
$taskSheet.Range("E5").Value = 0

$taskSheet.Range("A6").Value = "Project Build"
$taskSheet.Range("B6").Value = "Iteration Review with client"
$taskSheet.Range("C6").Value = 2
$taskSheet.Range("D6").Value = 1
$taskSheet.Range("E6").Value = 0

# --- ACTIVITY LOG SUMMARY SHEET ---
$summarySheet = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$summarySheet.Range("D1").Value = "Jesse Hare"

$summarySheet.Range("A4").Value = "Project Build"
$summarySheet.Range("B4").Value = 15
$summarySheet.Range("C4").Value = 23
